# Fruta / hortaliza, semanal
#
# The underlying data rows (2-21) on the sheet are re-shuffled: the tuple of
# (Fecha, Volumen, Precio minimo, Precio maximo, Precio promedio ponderado,
# Origen, Precio $/Kg) -> columns D, J, K, L, M, O, P -- moves from one row
# to another, while the rest of each row (Mercado ID, Mercado, Region,
# Codreg, Categoria ID, Categoria, Variedad, Calidad, Unidad de
# comercializacion, Kg o Unidades, Clasificacion) stays fixed, since those
# columns are constant for every data row anyway.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# after-row -> before-row: which row's (D,J,K,L,M,O,P) tuple ends up in
# this row once the edit is applied.
$map = @{}
$map[2]  = 15
$map[3]  = 21
$map[4]  = 10
$map[5]  = 7
$map[6]  = 6
$map[7]  = 12
$map[8]  = 18
$map[9]  = 16
$map[10] = 14
$map[11] = 13
$map[12] = 3
$map[13] = 4
$map[14] = 8
$map[15] = 19
$map[16] = 2
$map[17] = 5
$map[18] = 20
$map[19] = 17
$map[20] = 11
$map[21] = 9

# Snapshot the movable columns for every data row before mutating anything.
$snapshot = @{}
for ($r = 2; $r -le 21; $r++) {
    $snapshot[$r] = @{
        D = $ws.Cells.Item($r, 4).Value2
        J = $ws.Cells.Item($r, 10).Value2
        K = $ws.Cells.Item($r, 11).Value2
        L = $ws.Cells.Item($r, 12).Value2
        M = $ws.Cells.Item($r, 13).Value2
        O = $ws.Cells.Item($r, 15).Value2
        P = $ws.Cells.Item($r, 16).Value2
    }
}

# Write each row's new tuple from the snapshot of its mapped source row.
for ($r = 2; $r -le 21; $r++) {
    $src = $snapshot[$map[$r]]
    $ws.Cells.Item($r, 4).Value  = $src.D
    $ws.Cells.Item($r, 10).Value = $src.J
    $ws.Cells.Item($r, 11).Value = $src.K
    $ws.Cells.Item($r, 12).Value = $src.L
    $ws.Cells.Item($r, 13).Value = $src.M
    $ws.Cells.Item($r, 15).Value = $src.O
    $ws.Cells.Item($r, 16).Value = $src.P
}
